$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p016r_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p016r_1</id>", 2)
$d.Content.Find.Execute("<id>p016v_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p016v_1</id>", 2)
$d.Content.Find.Execute("<id>p016v_a2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p016v_2</id>", 2)
$d.Content.Find.Execute("<id>p016v_a3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p016v_3</id>", 2)
$d.Content.Find.Execute("<id>p016v_a4</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p016v_4</id>", 2)
$d.Content.Find.Execute("<id>p016v_a5</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p016v_5</id>", 2)
